# Refresh the cryptos list (GitHub Actions scheduled data pull).
# Updates the Price (D) / Volume(1h) (E) columns for each coin row, and
# reorders a couple of rows whose ranking swapped (Toncoin/WstETH,
# Cronos/Monero/PolygonEcosystemToken/EthereumClassic) by rewriting their
# Coin/Link/Price/Volume cells in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (losing trailing zeros / exact text form).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the refreshed cryptos feed.
$ws.Range("D2").Value = '76.763.98'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.953.39'
$ws.Range("E3").Value = '  +2.90%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '199.90'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").Value = '595.18'
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.549'
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D10").Value = '2.949.53'
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("D11").Value = '0.448'
$ws.Range("E11").Value = '  +14.31%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.493.56'
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = '4.90'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = '76.664.53'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").Value = '28.15'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '2.941.04'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("D19").Value = '13.34'
$ws.Range("E19").Value = '  +6.68%  '
$ws.Range("D20").Value = '8.67'
$ws.Range("E20").Value = '  -4.82%  '
$ws.Range("D21").Value = '371.10'
$ws.Range("E21").Value = '  -3.51%  '
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  +4.27%  '
$ws.Range("D23").Value = '2.26'
$ws.Range("E23").Value = '  -3.81%  '
$ws.Range("D24").Value = '72.40'
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '3.088.01'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("D27").Value = '4.25'
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").Value = '9.66'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  +6.21%  '
$ws.Range("D32").Value = '1.37'
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").Value = '494.33'
$ws.Range("E33").Value = '  -3.92%  '
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  +23.68%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '166.77'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '0.398'
$ws.Range("E38").Value = '  +15.29%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '20.15'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("E41").Value = '  -6.80%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '180.66'
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '4.91'
$ws.Range("E44").Value = '  -3.84%  '
$ws.Range("D45").Value = '1.64'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").Value = '40.11'
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("E47").Value = '  -4.70%  '
$ws.Range("D48").Value = '0.589'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("E49").Value = '  +3.59%  '
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").Value = '22.58'
$ws.Range("E51").Value = '  +4.30%  '
